$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.230.24'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").Value = '1.655.99'
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.31'
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5241'
$ws.Range("E6").Value = '  -2.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  -0.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2671'
$ws.Range("E8").Value = '  +0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06363'
$ws.Range("E9").Value = '  -0.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.69'
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07761'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.580'
$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("D13").Value = '1.655.68'
$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("D14").Value = '1.885.91'
$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5672'
$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").Value = '0.0₅8197'
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.52'
$ws.Range("E17").Value = '  -1.46%  '

$ws.Range("D18").Value = '26.246.86'

$ws.Range("E19").Value = '  -0.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.714'
$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.30'
$ws.Range("E21").Value = '  -3.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.034'
$ws.Range("E23").Value = '  -0.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.48'
$ws.Range("E25").Value = '  -2.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1205'
$ws.Range("E26").Value = '  -2.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.284'
$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.04'
$ws.Range("E28").Value = '  -1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.491'
$ws.Range("E29").Value = '  -0.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05630'
$ws.Range("E30").Value = '  -4.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.277'
$ws.Range("E31").Value = '  -0.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.499'
$ws.Range("E32").Value = '  -2.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.374'
$ws.Range("E33").Value = '  +1.78%  '

$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.806'
$ws.Range("E35").Value = '  -1.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9465'
$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5772'
$ws.Range("E38").Value = '  -1.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01600'
$ws.Range("E39").Value = '  -1.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.918'
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.580'
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8474'
$ws.Range("E42").Value = '  -2.28%  '

$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("D44").Value = '1.029.60'
$ws.Range("E44").Value = '  -4.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.29'
$ws.Range("E45").Value = '  -1.93%  '

$ws.Range("D46").Value = '1.796.55'
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '58.58'
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("E49").Value = '  -0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05318'
$ws.Range("E50").Value = '  +2.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.028'
$ws.Range("E51").Value = '  -0.43%  '
